# Scheduled runner update: refresh Marketboard-derived profit figures (columns H-N)
# across the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 415
$ws.Range("I4").Value = 420.33334
$ws.Range("J4").Value = 399
$ws.Range("K4").Value = 420.33334
$ws.Range("L4").Value = 399
$ws.Range("M4").Value = -306.33334
$ws.Range("N4").Value = -627
$ws.Range("H19").Value = 4902
$ws.Range("I19").Value = 4805.3076
$ws.Range("J19").Value = 5027.7
$ws.Range("K19").Value = 4805.3076
$ws.Range("L19").Value = 5027.7
$ws.Range("M19").Value = -4630.3076
$ws.Range("N19").Value = -5377.7
$ws.Range("H33").Value = 733.1429000000001
$ws.Range("J33").Value = 669.8333
$ws.Range("L33").Value = 669.8333
$ws.Range("N33").Value = -1127.8333
$ws.Range("H74").Value = 2061
$ws.Range("I74").Value = 2061
$ws.Range("K74").Value = 2061
$ws.Range("M74").Value = -1125
$ws.Range("H77").Value = 2061
$ws.Range("I77").Value = 2061
$ws.Range("K77").Value = 10305
$ws.Range("M77").Value = -5625
$ws.Range("H97").Value = 821.9
$ws.Range("J97").Value = 817.8889
$ws.Range("L97").Value = 2453.6667
$ws.Range("N97").Value = -3445.6667
$ws.Range("H111").Value = 1148.7142
$ws.Range("I111").Value = 999.4286
$ws.Range("K111").Value = 2998.2858
$ws.Range("M111").Value = 68.71420000000035
$ws.Range("H137").Value = 14762.292
$ws.Range("I137").Value = 26152.834
$ws.Range("K137").Value = 78458.50199999999
$ws.Range("M137").Value = -75908.50199999999
$ws.Range("H138").Value = 2193.2258
$ws.Range("J138").Value = 4319.091
$ws.Range("L138").Value = 12957.273
$ws.Range("N138").Value = -23237.273
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1487.8125
$ws.Range("I132").Value = 1052.48
$ws.Range("J132").Value = 3042.5715
$ws.Range("K132").Value = 3157.44
$ws.Range("L132").Value = 9127.7145
$ws.Range("M132").Value = -627.4400000000001
$ws.Range("N132").Value = -14187.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3884
$ws.Range("I105").Value = 2932.3333
$ws.Range("J105").Value = 4359.8335
$ws.Range("K105").Value = 2932.3333
$ws.Range("L105").Value = 4359.8335
$ws.Range("M105").Value = -1185.3333
$ws.Range("N105").Value = -7853.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3450392.8
$ws.Range("I31").Value = 4763586.5
$ws.Range("K31").Value = 4763586.5
$ws.Range("M31").Value = -4763291.5
$ws.Range("H34").Value = 3450392.8
$ws.Range("I34").Value = 4763586.5
$ws.Range("K34").Value = 4763586.5
$ws.Range("M34").Value = -4763384.5
$ws.Range("H58").Value = 1702.1072
$ws.Range("I58").Value = 1402.2916
$ws.Range("K58").Value = 1402.2916
$ws.Range("M58").Value = -1199.2916
$ws.Range("H63").Value = 35000
$ws.Range("J63").Value = 35000
$ws.Range("L63").Value = 35000
$ws.Range("N63").Value = -36372
$ws.Range("H66").Value = 35000
$ws.Range("J66").Value = 35000
$ws.Range("L66").Value = 105000
$ws.Range("N66").Value = -111864
$ws.Range("H68").Value = 35000
$ws.Range("J68").Value = 35000
$ws.Range("L68").Value = 35000
$ws.Range("N68").Value = -36498
$ws.Range("H69").Value = 59000
$ws.Range("I69").Value = 28000
$ws.Range("J69").Value = 90000
$ws.Range("K69").Value = 28000
$ws.Range("L69").Value = 90000
$ws.Range("M69").Value = -27251
$ws.Range("N69").Value = -91498
$ws.Range("H71").Value = 35000
$ws.Range("J71").Value = 35000
$ws.Range("L71").Value = 105000
$ws.Range("N71").Value = -112488
$ws.Range("H72").Value = 59000
$ws.Range("I72").Value = 28000
$ws.Range("J72").Value = 90000
$ws.Range("K72").Value = 84000
$ws.Range("L72").Value = 270000
$ws.Range("M72").Value = -80256
$ws.Range("N72").Value = -277488
$ws.Range("H74").Value = 79997
$ws.Range("J74").Value = 89996.39999999999
$ws.Range("L74").Value = 89996.39999999999
$ws.Range("N74").Value = -91744.39999999999
$ws.Range("H77").Value = 79997
$ws.Range("J77").Value = 89996.39999999999
$ws.Range("L77").Value = 269989.2
$ws.Range("N77").Value = -278725.2
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H122").Value = 1877.2354
$ws.Range("I122").Value = 1869.6875
$ws.Range("K122").Value = 5609.0625
$ws.Range("M122").Value = -3159.0625
$ws.Range("H132").Value = 18385.486
$ws.Range("I132").Value = 19893.812
$ws.Range("K132").Value = 59681.436
$ws.Range("M132").Value = -57151.436
$ws.Range("H134").Value = 1587.7632
$ws.Range("I134").Value = 1533.1562
$ws.Range("J134").Value = 1879
$ws.Range("K134").Value = 4599.4686
$ws.Range("L134").Value = 5637
$ws.Range("M134").Value = -2064.4686
$ws.Range("N134").Value = -10707
$ws.Range("H136").Value = 1702.1072
$ws.Range("I136").Value = 1402.2916
$ws.Range("K136").Value = 4206.8748
$ws.Range("M136").Value = -1656.8748
$ws.Range("H138").Value = 88415.06
$ws.Range("J138").Value = 88415.06
$ws.Range("L138").Value = 88415.06
$ws.Range("N138").Value = -98695.06
$ws.Range("H140").Value = 100104.664
$ws.Range("J140").Value = 100104.664
$ws.Range("L140").Value = 100104.664
$ws.Range("N140").Value = -110464.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 912.1
$ws.Range("I5").Value = 880.1111
$ws.Range("K5").Value = 2640.3333
$ws.Range("M5").Value = -2528.3333
$ws.Range("H117").Value = 854.44446
$ws.Range("J117").Value = 827.75
$ws.Range("L117").Value = 2483.25
$ws.Range("N117").Value = -9367.25
$ws.Range("H135").Value = 912.1
$ws.Range("I135").Value = 880.1111
$ws.Range("K135").Value = 7920.9999
$ws.Range("M135").Value = -5385.9999
$ws.Range("H138").Value = 9066.105
$ws.Range("I138").Value = 9765
$ws.Range("K138").Value = 29295
$ws.Range("M138").Value = -24155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5261.087
$ws.Range("I80").Value = 2583.8667
$ws.Range("J80").Value = 10280.875
$ws.Range("K80").Value = 2583.8667
$ws.Range("L80").Value = 10280.875
$ws.Range("M80").Value = -1585.8667
$ws.Range("N80").Value = -12276.875
$ws.Range("H83").Value = 5261.087
$ws.Range("I83").Value = 2583.8667
$ws.Range("J83").Value = 10280.875
$ws.Range("K83").Value = 12919.3335
$ws.Range("L83").Value = 51404.375
$ws.Range("M83").Value = -7927.333500000001
$ws.Range("N83").Value = -61388.375
$ws.Range("H122").Value = 4080.4119
$ws.Range("I122").Value = 3585.3809
$ws.Range("J122").Value = 4880.077
$ws.Range("K122").Value = 10756.1427
$ws.Range("L122").Value = 14640.231
$ws.Range("M122").Value = -8306.1427
$ws.Range("N122").Value = -19540.231
$ws.Range("H132").Value = 6184.273
$ws.Range("I132").Value = 6652.7
$ws.Range("K132").Value = 19958.1
$ws.Range("M132").Value = -17428.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4316.304
$ws.Range("J46").Value = 5423.6875
$ws.Range("L46").Value = 5423.6875
$ws.Range("N46").Value = -5799.6875
$ws.Range("H68").Value = 4163.304
$ws.Range("I68").Value = 3073
$ws.Range("K68").Value = 3073
$ws.Range("M68").Value = -2324
$ws.Range("H71").Value = 4163.304
$ws.Range("I71").Value = 3073
$ws.Range("K71").Value = 15365
$ws.Range("M71").Value = -11621
$ws.Range("H132").Value = 2156.625
$ws.Range("I132").Value = 2208.8333
$ws.Range("K132").Value = 6626.499899999999
$ws.Range("M132").Value = -4096.499899999999
$ws.Range("H135").Value = 110000
$ws.Range("J135").Value = 110000
$ws.Range("L135").Value = 110000
$ws.Range("N135").Value = -120140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 767.95
$ws.Range("J113").Value = 941.44446
$ws.Range("L113").Value = 2824.33338
$ws.Range("N113").Value = -7164.33338
$ws.Range("H122").Value = 38968.6
$ws.Range("I122").Value = 46233
$ws.Range("J122").Value = 3857.3333
$ws.Range("K122").Value = 138699
$ws.Range("L122").Value = 11571.9999
$ws.Range("M122").Value = -136249
$ws.Range("N122").Value = -16471.9999
$ws.Range("H132").Value = 3542.7742
$ws.Range("I132").Value = 3586.2222
$ws.Range("K132").Value = 10758.6666
$ws.Range("M132").Value = -8228.6666
